# Update "想去人数" (number of people interested) figures to the newly
# scraped values for gh-pages output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1337
$wsExhibit.Range("F5").Value = 99
$wsExhibit.Range("F8").Value = 11578
$wsExhibit.Range("F9").Value = 4365
$wsExhibit.Range("F11").Value = 33
$wsExhibit.Range("F14").Value = 2540
$wsExhibit.Range("F15").Value = 1088
$wsExhibit.Range("F16").Value = 136
$wsExhibit.Range("F17").Value = 31
$wsExhibit.Range("F18").Value = 4377
$wsExhibit.Range("F21").Value = 11320
$wsExhibit.Range("F22").Value = 11239

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 1

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1337
$wsAll.Range("F5").Value = 99
$wsAll.Range("F8").Value = 11578
$wsAll.Range("F9").Value = 4365
$wsAll.Range("F11").Value = 33
$wsAll.Range("F14").Value = 2540
$wsAll.Range("F15").Value = 1
$wsAll.Range("F16").Value = 1088
$wsAll.Range("F17").Value = 136
$wsAll.Range("F18").Value = 31
$wsAll.Range("F19").Value = 4377
$wsAll.Range("F22").Value = 11320
$wsAll.Range("F23").Value = 11239

$wb.Save()
